$d = $word.ActiveDocument

# Translation update: "RPC Explorer" -> "Insight Explorer"
$d.Content.Find.Execute("RPC Explorer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Insight Explorer", 2)
